$d = $word.ActiveDocument

# Locate the exact span to replace: "for our brand-new article entitiled:"
# (covers the misspelled word and its surrounding proofErr markers) and
# rewrite it as separate runs, inserting the new "than look no further
# than" sentence and fixing "entitiled" -> "entitled".
$r = $d.Content
$found = $r.Find.Execute("for our brand-new article entitiled:")
if (-not $found) {
    throw "Could not find target text to replace"
}

# Re-wrap the found span in a fresh Range object (anchored at the same
# Start/End) before calling InsertXML - calling InsertXML directly on the
# Range returned by Find.Execute does not replace in place.
$target = $d.Range($r.Start, $r.End)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
       '<w:r><w:t>than look no further than</w:t></w:r>' + `
       '<w:r><w:t xml:space="preserve"> our brand-new article </w:t></w:r>' + `
       '<w:r><w:t>entitled</w:t></w:r>' + `
       '<w:r><w:t>:</w:t></w:r>' + `
       '</w:p>'

$target.InsertXML($xml)
